# Applies the content updates described in the diff for
# sequencing_spreadsheet_template.by_provider.xlsx ("Examples & Info" sheet).
#
# Row 2 holds column descriptions; several cells previously contained the
# placeholder text "another description" and are now filled in with the
# real descriptive text. AR2 also has a typo fix. Row 6 gains a few
# regex-validation strings for date columns that were previously empty.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Examples & Info")

# --- Row 2: column descriptions ---
$ws.Range("L2").Value = "Location of the library/experiement description"
$ws.Range("M2").Value = "Plate label used"
$ws.Range("T2").Value = "Any treatment / pertubation applied to the individual/ cell-line."
$ws.Range("Y2").Value = "How was the mouse/animal killed?"
$ws.Range("AH2").Value = "Well position of the barcode"
$ws.Range("AJ2").Value = "Well position of the barcode"
$ws.Range("AO2").Value = "Nucleic acid concentration"
$ws.Range("AP2").Value = "Library molarity"
$ws.Range("AR2").Value = "Average fragment length (basepairs)"

# --- Row 6: validation regex examples for date columns ---
$ws.Range("W6").Value = "[0-9]{4}-[0-9]{2}-[0-9]{2}"
$ws.Range("X6").Value = "[0-9]{4}-[0-9]{2}-[0-9]{2}"
$ws.Range("AT6").Value = "[0-9]{4}-[0-9]{2}-[0-9]{2}"

$wb.Save()
